$d = $word.ActiveDocument
$d.Content.Find.Execute("(BNCC - EF02MA22)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "(EF02MA22)", 2)
